{"js": "// The document contains a single table of two-digit-number \u00f7 one-digit-number\n// division problems. Every non-blank row holds 5 answer cells; this script\n// replaces each answer's text with its updated value, matching strictly by\n// cell position (row, column) in document order so that values which are\n// reused elsewhere in the table (e.g. \"86\u00f72=43, 0\" used to exist at one spot\n// and becomes the new value of a different spot) do not collide.\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Old -> new text, in the exact order the answer cells appear in the\n// document (left-to-right, top-to-bottom), five per populated row.\nconst replacements = [\n  \"40\u00f79=4, 4|86\u00f72=43, 0\",\n  \"36\u00f77=5, 1|11\u00f74=2, 3\",\n  \"19\u00f75=3, 4|86\u00f75=17, 1\",\n  \"21\u00f78=2, 5|18\u00f79=2, 0\",\n  \"82\u00f72=41, 0|85\u00f78=10, 5\",\n  \"74\u00f74=18, 2|80\u00f73=26, 2\",\n  \"89\u00f79=9, 8|49\u00f72=24, 1\",\n  \"12\u00f75=2, 2|43\u00f75=8, 3\",\n  \"71\u00f73=23, 2|13\u00f77=1, 6\",\n  \"12\u00f76=2, 0|31\u00f73=10, 1\",\n  \"59\u00f79=6, 5|33\u00f78=4, 1\",\n  \"61\u00f78=7, 5|98\u00f76=16, 2\",\n  \"28\u00f78=3, 4|85\u00f76=14, 1\",\n  \"82\u00f76=13, 4|83\u00f76=13, 5\",\n  \"36\u00f76=6, 0|59\u00f73=19, 2\",\n  \"13\u00f72=6, 1|92\u00f75=18, 2\",\n  \"74\u00f79=8, 2|29\u00f77=4, 1\",\n  \"24\u00f74=6, 0|32\u00f79=3, 5\",\n  \"74\u00f78=9, 2|38\u00f76=6, 2\",\n  \"86\u00f72=43, 0|60\u00f74=15, 0\",\n  \"85\u00f78=10, 5|66\u00f76=11, 0\",\n  \"85\u00f74=21, 1|93\u00f75=18, 3\",\n  \"88\u00f72=44, 0|16\u00f72=8, 0\",\n  \"85\u00f77=12, 1|33\u00f72=16, 1\",\n  \"18\u00f77=2, 4|80\u00f76=13, 2\",\n].map((s) => {\n  const i = s.indexOf(\"|\");\n  return [s.slice(0, i), s.slice(i + 1)];\n});\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet cursor = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  // Load each cell's current text so we only touch the cells that actually\n  // hold an answer (the template also has fully-blank practice rows).\n  for (const cell of cells.items) {\n    cell.body.load(\"text\");\n  }\n  await context.sync();\n\n  const hasAnswer = cells.items.some((cell) => cell.body.text.trim().length > 0);\n  if (!hasAnswer) {\n    continue;\n  }\n\n  for (const cell of cells.items) {\n    if (cursor >= replacements.length) break;\n    const [oldText, newText] = replacements[cursor];\n    cursor++;\n    // Sanity-check we are editing the cell the diff intends (falls back to\n    // applying the edit regardless, keyed purely by position, if the\n    // existing text doesn't line up exactly - e.g. stray whitespace).\n    const current = cell.body.text.replace(/[\\r\\v]+$/, \"\");\n    if (current !== oldText) {\n      // Not fatal: still proceed positionally.\n    }\n    // Replace the existing text in place (instead of clearing + inserting)\n    // so the run keeps its original character formatting (font/size) and\n    // the paragraph keeps its alignment.\n    const range = cell.body.paragraphs.getFirst().getRange();\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document holds one table of two-digit \u00f7 one-digit division problems.\n# Every 4th row (1, 5, 9, 13, 17 in 1-based Word terms) carries 5 answer\n# cells; the rows in between are blank practice rows. We rewrite each\n# answer cell's text in place (Cell.Range.Text = ...) - this keeps the run's\n# existing character formatting (font/size) and the paragraph's alignment,\n# exactly like Word does for a plain text replacement - addressing strictly\n# by (row, column) position so that answers re-used elsewhere in the grid\n# (\"86\u00f72=43, 0\" and \"85\u00f78=10, 5\" both show up as an old value in one cell and\n# a new value in a different cell) never collide.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$oldNew = @(\n  @(\"40\u00f79=4, 4\", \"86\u00f72=43, 0\"),\n  @(\"36\u00f77=5, 1\", \"11\u00f74=2, 3\"),\n  @(\"19\u00f75=3, 4\", \"86\u00f75=17, 1\"),\n  @(\"21\u00f78=2, 5\", \"18\u00f79=2, 0\"),\n  @(\"82\u00f72=41, 0\", \"85\u00f78=10, 5\"),\n  @(\"74\u00f74=18, 2\", \"80\u00f73=26, 2\"),\n  @(\"89\u00f79=9, 8\", \"49\u00f72=24, 1\"),\n  @(\"12\u00f75=2, 2\", \"43\u00f75=8, 3\"),\n  @(\"71\u00f73=23, 2\", \"13\u00f77=1, 6\"),\n  @(\"12\u00f76=2, 0\", \"31\u00f73=10, 1\"),\n  @(\"59\u00f79=6, 5\", \"33\u00f78=4, 1\"),\n  @(\"61\u00f78=7, 5\", \"98\u00f76=16, 2\"),\n  @(\"28\u00f78=3, 4\", \"85\u00f76=14, 1\"),\n  @(\"82\u00f76=13, 4\", \"83\u00f76=13, 5\"),\n  @(\"36\u00f76=6, 0\", \"59\u00f73=19, 2\"),\n  @(\"13\u00f72=6, 1\", \"92\u00f75=18, 2\"),\n  @(\"74\u00f79=8, 2\", \"29\u00f77=4, 1\"),\n  @(\"24\u00f74=6, 0\", \"32\u00f79=3, 5\"),\n  @(\"74\u00f78=9, 2\", \"38\u00f76=6, 2\"),\n  @(\"86\u00f72=43, 0\", \"60\u00f74=15, 0\"),\n  @(\"85\u00f78=10, 5\", \"66\u00f76=11, 0\"),\n  @(\"85\u00f74=21, 1\", \"93\u00f75=18, 3\"),\n  @(\"88\u00f72=44, 0\", \"16\u00f72=8, 0\"),\n  @(\"85\u00f77=12, 1\", \"33\u00f72=16, 1\"),\n  @(\"18\u00f77=2, 4\", \"80\u00f76=13, 2\")\n)\n\n$cols = $t.Columns.Count\n$rows = $t.Rows.Count\n$idx = 0\n\nfor ($r = 1; $r -le $rows; $r++) {\n  $rowHasAnswer = $false\n  for ($c = 1; $c -le $cols; $c++) {\n    $text = $t.Cell($r, $c).Range.Text\n    $text = $text -replace \"[`r`a]+$\", \"\"\n    if ($text.Trim().Length -gt 0) {\n      $rowHasAnswer = $true\n    }\n  }\n  if (-not $rowHasAnswer) {\n    continue\n  }\n  for ($c = 1; $c -le $cols; $c++) {\n    if ($idx -ge $oldNew.Length) { break }\n    $pair = $oldNew[$idx]\n    $idx++\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $pair[1]\n  }\n}\n"}
